# Fixing rules with LocalDate problem
# The "Set Queue Enter Date" rows used toDate(java.time.LocalDate.now()) as the
# action's parameter, but setQueueEnterDate expects a java.time.LocalDate, not a
# java.util.Date. Remove the unnecessary toDate(...) wrapper so it passes the
# LocalDate straight through.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldValue = "setQueueEnterDate, toDate(java.time.LocalDate.now())"
$newValue = "setQueueEnterDate, java.time.LocalDate.now()"

# Column E holds the "Set Field Value" action for each queue-enter-date rule
# (rows 27-31). Update every cell that still has the old expression.
$usedRange = $ws.UsedRange
foreach ($cell in $usedRange.Cells) {
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
